$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Extend the "PickUpBehavior ... ArrayList." paragraph with the new
#    sentence about behavior priority.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "a new behavior class named PickUpBehavior will be added to the Zombie class behaviors ArrayList.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "a new behavior class named PickUpBehavior will be added to the Zombie class behaviors ArrayList. This behavior will be lower priority than AttackBehavior, but higher than HuntBehavior or WanderBehavior.",
    2)

# ---------------------------------------------------------------------------
# 2. Extend the final paragraph ("... must then be checked.") with the two
#    new sentences describing how the ground check is implemented, using a
#    blank line (two manual line breaks) between them.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertAfter(
    " This can be done by using map.LocationOf(zombie) to get the zombie" + [char]8217 + "s current map location, then calling getItems() on that location to fetch all items on the ground at that location, then calling item.asWeapon() to check the item" + [char]8217 + "s status as a weapon." +
    [char]11 + [char]11 +
    "If there is a valid Weapon on the ground in the same location as the zombie, PickUpBehavior will return a pickUpItemAction for that Weapon, which Zombie will then execute."
)

# ---------------------------------------------------------------------------
# 3. Append the new sections at the end of the document:
#      (blank)
#      Implementing Zombie Moan
#      Every time the  playTurn method in Zombie is called, ...
#      (blank, centered)
#      BEATING UP THE ZOMBIES (centered)
#      (blank)
#      Implementing Zombie dismemberment
#      (blank)
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = [char]9 + "Implementing Zombie Moan"

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "Every time the  playTurn method in Zombie is called, before the program loops over the Behavior ArrayList, a random probability will return a 10% chance that the system prints a zombie moan."

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Alignment = 1

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Alignment = 1
$p.Range.Text = "BEATING UP THE ZOMBIES"

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = [char]9 + "Implementing Zombie dismemberment"

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)

Write-Output "done"
